$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07496433333333334
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 0.1720421704636667
$ws.Range("R2").Value = 1.548379534173
$ws.Range("S2").Value = 0.0158275801650097
$ws.Range("T2").Value = 0.0158275801650097

$ws.Range("G3").Value = 0.07496433333333334
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 8.365398195652888
$ws.Range("S3").Value = 0.769602070219672
$ws.Range("T3").Value = 0.7696020702196722

$ws.Range("G4").Value = 0.07496433333333334
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("S4").Value = 0.2145703496153182
$ws.Range("T4").Value = 0.2145703496153182
